# Bancolombia e-prepago data-driven workbook update:
#  - Add a new scenario row (row 3) to the "Datos" sheet, mirroring row 2 but
#    with a recharge value ("valorRecarga") above the maximum allowed amount.
#  - Add a new, empty worksheet "CargaSobreMaximo" after "Datos".

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Datos")

# --- 1. Populate the new data row (row 3) -----------------------------------
# Set the cell values first so purely-numeric entries (valorRecarga) are
# stored as numbers, matching the style of row 2.
$ws1.Range("A3").Value = 1
$ws1.Range("B3").Value = 93221450
$ws1.Range("C3").Value = 1
$ws1.Range("D3").Value = "autotest25"
$ws1.Range("E3").Value = 1234
$ws1.Range("F3").Value = 4321
$ws1.Range("G3").Value = "Acierto"
$ws1.Range("H3").Value = "000"
$ws1.Range("I3").Value = "0369"
$ws1.Range("J3").Value = "NO ERROR"
$ws1.Range("K3").Value = "bolp"
$ws1.Range("L3").Value = "ACTIVO"
$ws1.Range("M3").Value = 3000000
$ws1.Range("N3").Value = "Recargar"
$ws1.Range("O3").Value = "Ahorros"
$ws1.Range("P3").Value = "406-739740-05"

# Copy the formatting (styles) of row 2 onto row 3 so the new row looks the
# same as the existing data row.
$ws1.Range("A2:P2").Copy()
$ws1.Range("A3:P3").PasteSpecial(-4122) # xlPasteFormats
$excel.CutCopyMode = 0

# Match row 2's height.
$ws1.Rows.Item(3).RowHeight = $ws1.Rows.Item(2).RowHeight

# codigoError/codigoTransaccion ("000"/"0369") must stay text (leading
# zeros); re-apply after the format paste so they don't get reinterpreted as
# numbers.
$ws1.Range("H3").Value = "000"
$ws1.Range("I3").Value = "0369"

# --- 2. Add the new "CargaSobreMaximo" worksheet ----------------------------
$newSheet = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $ws1)
$newSheet.Name = "CargaSobreMaximo"

# --- 3. Re-activate "Datos" and restore its selection -----------------------
$ws1.Activate()
$ws1.Range("A1:Q1").Select()
